# Apply the "Add files via upload" edit to meta/en/1-1-1.xlsx
#
# Summary of the change (from the supplied OOXML diff):
#  1. The workbook window view (bookViews/workbookView) size/position changed.
#  2. The Data reporter contact block (rows 6-10, column B) was updated with a
#     new organization / contact person / email / phone / website.
#  3. The sheet's active selection moved from B24 to B6.
#
# The sheet is protected, but the data-entry cells in column B (rows 6-10)
# are explicitly unlocked in the original workbook, so they can be edited
# without needing to unprotect/re-protect the sheet (doing so would also
# risk altering the <sheetProtection .../> element that must stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2. Data reporter -------------------------------------------------
$ws.Range("B6").Value  = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value  = "Kalymbetova Yryskan"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# --- Workbook window geometry ------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.Left   = 0
    $win.Top    = 0
    $win.Width  = 1440
    $win.Height = 591.75
} catch {
    # Window geometry may not be adjustable in every host; ignore if so.
}

# --- Active cell / selection --------------------------------------------
$ws.Activate()
$ws.Range("B6").Select() | Out-Null
